# ----------------------------------------------------------------------------
# Excel COM-interop edit script
# Applies the "Okno serwisowe 23.04.2024" data-refresh update to
# global_sdg_indicators.xlsx: new 2023 (col S) figures, revised 2022 (col R)
# figures, two label tweaks, and small view/print-setup tweaks.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SDGs_1-17")

# --- 1. Cells that were previously empty and now carry a new value.
#     Each new cell must inherit the numeric-format style already used by the
#     rest of the row (column F always carries it), so we copy *formats only*
#     from the row-anchor cell before writing the value - this reuses the
#     existing style index instead of minting a new one.
$formatCopyPairs = @(
    @("F17", "S17", 15.2),
    @("F29", "R29", 61),
    @("F30", "R30", 12),
    @("F31", "R31", 10.5),
    @("F32", "R32", 38.4),
    @("F49", "R49", 0.82),
    @("F51", "R51", 2),
    @("F61", "S61", 46),
    @("F62", "S62", 13),
    @("F63", "S63", 33),
    @("F65", "S65", 140),
    @("F66", "S66", 68),
    @("F67", "S67", 72),
    @("F68", "R68", 426.2),
    @("F69", "R69", 28.8),
    @("F70", "R70", 253.9),
    @("F71", "R71", 24.6),
    @("F72", "R72", 11.9),
    @("F73", "R73", 20.5),
    @("F74", "R74", 3.4),
    @("F86", "R86", 1.2),
    @("F177", "R177", 46.1),
    @("F178", "R178", 33.1),
    @("F179", "R179", 42.4),
    @("F180", "R180", 32.299999999999997),
    @("F181", "R181", 56.3),
    @("F182", "R182", 35.4),
    @("F183", "R183", 59.7),
    @("F211", "R211", 16.88),
    @("F216", "R216", 2),
    @("F217", "R217", 2.5),
    @("F218", "R218", 1.4),
    @("F222", "R222", 42.3),
    @("F223", "R223", 44.62),
    @("F224", "R224", 39.799999999999997),
    @("F225", "R225", 74.430000000000007),
    @("F226", "R226", 55.29),
    @("F227", "R227", 40.81),
    @("F228", "R228", 32.590000000000003),
    @("F229", "R229", 26.02),
    @("F230", "R230", 28.35),
    @("F231", "R231", 31.77),
    @("F232", "R232", 34.49),
    @("F233", "R233", 25.38),
    @("F234", "R234", 29.22),
    @("F235", "R235", 39.78),
    @("F236", "R236", 45.09),
    @("F237", "R237", 44.47),
    @("F238", "R238", 40.81),
    @("F239", "R239", 40.99),
    @("F240", "R240", 57.5),
    @("F241", "S241", 2.8),
    @("F242", "S242", 2.8),
    @("F243", "S243", 2.9),
    @("F244", "S244", 11.3),
    @("F245", "S245", 3),
    @("F246", "S246", 2.1),
    @("F247", "S247", 2),
    @("F248", "S248", 1.7),
    @("F249", "S249", 5.2),
    @("F280", "R280", 9.7799999999999994),
    @("F319", "S319", 4.0599999999999996),
    @("F320", "S320", 2603.83),
    @("F321", "S321", 1744.48),
    @("F322", "S322", 0.38),
    @("F323", "S323", 6.36),
    @("F324", "S324", 254.4),
    @("F325", "S325", 62.44),
    @("F418", "R418", 0.49),
    @("F463", "S463", 2603.8000000000002),
    @("F464", "S464", 1.26),
    @("F469", "S469", 49.6),
    @("F470", "S470", -5.0999999999999996),
)

foreach ($pair in $formatCopyPairs) {
    $srcRef, $dstRef, $newVal = $pair
    $ws.Range($srcRef).Copy() | Out-Null
    $ws.Range($dstRef).PasteSpecial(-4122) | Out-Null
    $ws.Range($dstRef).Value = $newVal
}
$excel.CutCopyMode = $false

# --- 2. Cells that already had a value and are simply being revised.
$simpleValueUpdates = @(
    @("Q49", 0.46),
    @("R65", 69),
    @("R66", 36),
    @("Q211", 15.61),
    @("R214", 6.1),
    @("R215", 10.6),
    @("M216", 5.4),
    @("M217", 5.9),
    @("M218", 4.5999999999999996),
    @("Q242", 3.4),
    @("Q244", 12),
    @("Q245", 4),
    @("R245", 3.2),
    @("Q246", 2.6),
    @("Q249", 6.1),
    @("R249", 5.4),
    @("R278", 17.8),
    @("N319", 3.29),
    @("Q319", 3.98),
    @("R319", 4.1100000000000003),
    @("N464", 1.51),
    @("P464", 1.46),
    @("R464", 1.28),
    @("R469", 49.2),
    @("R470", -3.4),
)

foreach ($pair in $simpleValueUpdates) {
    $ref, $newVal = $pair
    $ws.Range($ref).Value = $newVal
}

# --- 3. Two shared-string label updates.
$ws.Range("B471").Value = "17.14.1 Number of countries with mechanisms in place to enhance policy coherence of sustainable development - PROXY!"
$ws.Range("A479").Value = "Last update: 23-04-2024, 13:07"

# --- 4. View state: selection moves from A4 to the A1:B1 header range.
$ws.Activate()
$ws.Range("A1:B1").Select()

# --- 5. Print setup: scale down from 74% to 72%.
$ws.PageSetup.Zoom = 72

Write-Host "edit.ps1 completed"
